$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("N2").Value = 1.26
$ws.Range("P2").Value = 1.25
$ws.Range("S2").Value = 1.32

$ws.Range("P3").Value = 1.96

$ws.Range("Q6").Value = 3

$ws.Range("G8").Value = 2.42
$ws.Range("I8").Value = 4.1

$ws.Range("G9").Value = 2.2
$ws.Range("P9").Value = 1.7
$ws.Range("Q9").Value = 2.2
